# Update "gh-pages" data workbook: refresh counters for existing events,
# and insert the new "庐江·夏日游嘉年华" event (2024-07-27) as the new
# row 21 in both the "展览" sheet and the "全部类型" sheet, shifting the
# rows below it down by one and renumbering the running index in column A.

$wb = $excel.ActiveWorkbook

function Update-Sheet($ws, $lastRowBeforeInsert) {

    # ---- 1. refresh "want to go" counters (column F) for the untouched rows ----
    $counterUpdates = @{
        2  = 1601
        3  = 8970
        4  = 102
        6  = 684
        9  = 45
        10 = 75
        11 = 3820
        12 = 59
        13 = 379
        15 = 4232
        17 = 157
        19 = 8
    }
    foreach ($row in $counterUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $counterUpdates[$row]
    }

    # ---- 2. insert a new row at position 21, pushing every later row down ----
    $ws.Rows.Item(21).Insert()

    # ---- 3. populate the newly inserted row 21 with the new event ----
    $ws.Cells.Item(21, 1).Value = 20
    # match the bold / centered / bordered look used by every other index cell in column A
    $ws.Cells.Item(21, 1).Font.Bold = $true
    $ws.Cells.Item(21, 1).HorizontalAlignment = -4108
    $ws.Cells.Item(21, 1).VerticalAlignment = -4160
    $ws.Cells.Item(21, 1).Borders.LineStyle = 1

    $ws.Cells.Item(21, 2).NumberFormat = "@"
    $ws.Cells.Item(21, 2).Value = "2024-07-27"
    $ws.Cells.Item(21, 2).Style = "Normal"

    $ws.Cells.Item(21, 3).Value = "庐江·夏日游嘉年华"
    $ws.Cells.Item(21, 4).Value = "白山路东150米 庐江体育馆"

    $ws.Cells.Item(21, 5).NumberFormat = "@"
    $ws.Cells.Item(21, 5).Value = "2024.07.27 09:00-07.28 17:00"
    $ws.Cells.Item(21, 5).Style = "Normal"

    $ws.Cells.Item(21, 6).Value = 0
    $ws.Cells.Item(21, 7).Value = 60
    $ws.Cells.Item(21, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87569"
    $ws.Cells.Item(21, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/5tB3RWrN1718243791381.jpeg"

    # ---- 4. renumber column A (running index) for every row pushed down ----
    $lastRowAfterInsert = $lastRowBeforeInsert + 1
    for ($row = 22; $row -le $lastRowAfterInsert; $row++) {
        $ws.Cells.Item($row, 1).Value = $row - 1
    }

    # ---- 5. apply the counter bumps that land on the shifted rows ----
    # row that now holds "合肥·咒术回战only" (want-to-go 239 -> 240)
    $ws.Cells.Item(22, 6).Value = 240

    # row that now holds "合肥·第七届环形宇宙动漫游戏嘉年华" (want-to-go 2618 -> 2623)
    $ws.Cells.Item(24, 6).Value = 2623

    # row that now holds "合肥·银魂主题派对only2.0" (want-to-go 105 -> 109) -- always the last row
    $ws.Cells.Item($lastRowAfterInsert, 6).Value = 109
}

# "展览" sheet (sheet1.xml): rows populated 1..24 before the insert
$wsExhibit = $wb.Worksheets.Item(1)
Update-Sheet $wsExhibit 24

# "全部类型" sheet (sheet4.xml): rows populated 1..25 before the insert
$wsAll = $wb.Worksheets.Item(4)
Update-Sheet $wsAll 25
